$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.305.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.229.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.37"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.80%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.18"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.14"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.46"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.852"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.247.33"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.121.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +13.02%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.13"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +34.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.34"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.79"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.40%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.04"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.09"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +16.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0806"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.52%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.52"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0307"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.53%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.15%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.00%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.57%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.993"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.75%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.82%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.87%  "
